$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) was populated with a mangled literal string
# "5-5-2011-12" (folder-name leftover) instead of the actual game date.
# Correct it to the real ISO date "2012-05-05" for every data row.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    if ($cell.Value2 -eq "5-5-2011-12") {
        # Prefix with an apostrophe so Excel stores the corrected date as
        # literal text ("2012-05-05") instead of auto-converting it to a
        # date serial number, then strip the transient "quote prefix"
        # formatting it leaves behind so the cell's style is untouched.
        $cell.Value2 = "'2012-05-05"
        $cell.ClearFormats()
    }
}
